$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Section: "Current Release" notes block.
# Before:
#   7  Current Release (Heading2)
#   8  Removed unused Cinch.WPF MediatorSingleton ... (bookmark _GoBack)
#   9  <empty>
#   10 66413 Release : 27/11/11 5:53PM Uk time (Heading2)
#   11 Fixed small problem in TabControlEx ...
#
# After:
#   7  Current Release (Heading2)                       [unchanged]
#   8  Added some WPF message box enhancements ...       [new]
#   9  In the Cinch users (Michel Renaud)  own words      [new]
#   10 Early on when we started using Cinch ... (italic, bookmark _GoBack)
#   11 <empty>
#   12 66567 Release : 04/12/11 10:29AM Uk time (Heading2)
#   13 Removed unused Cinch.WPF MediatorSingleton ...     [relocated]
#   14 <empty>                                            [new]
#   15 66413 Release : 27/11/11 5:53PM Uk time (Heading2) [relocated]
#   16 Fixed small problem in TabControlEx ...            [unchanged]
# ------------------------------------------------------------------

# 1. Insert the two new plain paragraphs before the "Removed unused" paragraph.
$pRemoved = $d.Paragraphs(8)
$ins1 = $d.Range($pRemoved.Range.Start, $pRemoved.Range.Start)
$ins1.InsertBefore("In the Cinch users (Michel Renaud)  own words`r")
$ins2 = $d.Range($pRemoved.Range.Start, $pRemoved.Range.Start)
$ins2.InsertBefore("Added some WPF message box enhancements from a Cinch user, that seems quite cool.`r")

# 2. The old "Removed unused..." paragraph (still holding bookmark _GoBack) is now
#    paragraph 10 - turn its content into the new italic commentary text, keeping
#    the bookmark (and paragraph) exactly where it is.
$pItalic = $d.Paragraphs(10)
$italicRange = $d.Range($pItalic.Range.Start, $pItalic.Range.End)
$italicRange.Text = 'Early on when we started using Cinch, we noticed some limitations in the message box service, namely hardcoded captions (we sometimes need something more explicit, and the apps are bilingual too) and inability to set a specific button as the default button (e.g. having "Cancel" as the default instead of "Ok" ). I ended up adding more methods to the service to take care of that.'
$italicRange2 = $d.Range($pItalic.Range.Start, $pItalic.Range.End)
$italicRange2.Font.Italic = 1

# 3. Update the heading date/build-number paragraph (still paragraph 12) to the
#    new release heading.
$pHeading = $d.Paragraphs(12)
$headRange = $d.Range($pHeading.Range.Start, $pHeading.Range.End)
$headRange.Text = "66567 Release : 04/12/11 10:29AM Uk time"

# 4. Re-insert the original "Removed unused..." paragraph text right after the
#    new heading (it now precedes paragraph 13, "Fixed small problem...").
$pFixed = $d.Paragraphs(13)
$ins3 = $d.Range($pFixed.Range.Start, $pFixed.Range.Start)
$ins3.InsertBefore("Removed unused Cinch.WPF MediatorSingleton that was npt part of solution but existing on disk. My bad, Cinch.WPF always uses common linked MediatorSingleton file. This does not effect the Dlls just the files on disk, so no harm done here.`r")

# 5. Re-insert the original heading ("66413 Release : ...") plus the blank
#    paragraph that used to sit above it, again right before "Fixed small
#    problem...".
$pFixed2 = $d.Paragraphs(14)
$ins4 = $d.Range($pFixed2.Range.Start, $pFixed2.Range.Start)
$ins4.InsertBefore("66413 Release : 27/11/11 5:53PM Uk time`r")
$ins5 = $d.Range($pFixed2.Range.Start, $pFixed2.Range.Start)
$ins5.InsertBefore("`r")

# 6. Fix up styles/formatting of the relocated & inserted paragraphs - they
#    inherited "List Paragraph" bullet formatting from the following bullet
#    item, which is not what we want.
$pRemovedMoved = $d.Paragraphs(13)
$pRemovedMoved.Style = "Normal"
$pRemovedMoved.Range.ListFormat.RemoveNumbers()

$pBlank = $d.Paragraphs(14)
$pBlank.Style = "Normal"
$pBlank.Range.ListFormat.RemoveNumbers()

$pHeading2 = $d.Paragraphs(15)
$pHeading2.Style = "Heading 2"
$pHeading2.Range.ListFormat.RemoveNumbers()
$pHeading2.Range.Font.Name = "Tahoma"
$pHeading2.Range.Font.NameAscii = "Tahoma"
$pHeading2.Range.Font.NameOther = "Tahoma"
